$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 4544
$ws.Range("L3").Value = 4876
$ws.Range("B4").Value = 1715
$ws.Range("D4").Value = 1995
$ws.Range("L4").Value = 1206
$ws.Range("L5").Value = 282
$ws.Range("L6").Value = 4160
$ws.Range("B7").Value = 23347
$ws.Range("D7").Value = 28186
$ws.Range("L7").Value = 15068

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L3").Value = 338
$ws.Range("L7").Value = 1008

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L4").Value = 15
$ws.Range("L7").Value = 329

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L3").Value = 233
$ws.Range("L6").Value = 212
$ws.Range("L7").Value = 685

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 166
$ws.Range("L3").Value = 192
$ws.Range("L6").Value = 151
$ws.Range("L7").Value = 557

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L2").Value = 102
$ws.Range("L3").Value = 91
$ws.Range("L6").Value = 79
$ws.Range("L7").Value = 289

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L6").Value = 57
$ws.Range("L7").Value = 258

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 70
$ws.Range("L2").Value = 21

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L7").Value = 497
$ws.Range("L8").Value = 1008
$ws.Range("L11").Value = 241
$ws.Range("L12").Value = 35
$ws.Range("L14").Value = 74
$ws.Range("L15").Value = 113
$ws.Range("L23").Value = 163
$ws.Range("L25").Value = 87
$ws.Range("L26").Value = 15
$ws.Range("E29").Value = 1673
$ws.Range("L29").Value = 836
$ws.Range("L30").Value = 70
$ws.Range("L33").Value = 685
$ws.Range("L36").Value = 193
$ws.Range("L37").Value = 557
$ws.Range("L40").Value = 40
$ws.Range("L41").Value = 69
$ws.Range("L42").Value = 490
$ws.Range("L44").Value = 109
$ws.Range("L50").Value = 76
$ws.Range("L52").Value = 308
$ws.Range("B63").Value = 419
$ws.Range("D63").Value = 375
$ws.Range("E63").Value = 389
$ws.Range("L65").Value = 289
$ws.Range("L67").Value = 517
$ws.Range("L68").Value = 48
$ws.Range("L73").Value = 120
$ws.Range("L76").Value = 235
$ws.Range("L79").Value = 398
$ws.Range("L83").Value = 329
$ws.Range("L85").Value = 773
$ws.Range("L92").Value = 43
$ws.Range("L93").Value = 79
$ws.Range("L96").Value = 168
$ws.Range("L97").Value = 130
$ws.Range("L99").Value = 258
$ws.Range("B101").Value = 23347
$ws.Range("D101").Value = 28186
$ws.Range("L101").Value = 15068

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L2").Value = 150
$ws.Range("L7").Value = 517

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 251
$ws.Range("E4").Value = 77
$ws.Range("L4").Value = 38
$ws.Range("L6").Value = 219
$ws.Range("E7").Value = 1673
$ws.Range("L7").Value = 836

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("L3").Value = 31
$ws.Range("L7").Value = 109
$ws.Range("L2").Value = 45

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L7").Value = 235
$ws.Range("L3").Value = 43

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("L2").Value = 33
$ws.Range("L7").Value = 87

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("L7").Value = 74

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("L3").Value = 24
$ws.Range("L7").Value = 69

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L3").Value = 167
$ws.Range("L4").Value = 35
$ws.Range("L6").Value = 135
$ws.Range("L7").Value = 490

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L2").Value = 41
$ws.Range("L4").Value = 16
$ws.Range("L6").Value = 45
$ws.Range("L7").Value = 163

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("L6").Value = 48
$ws.Range("L7").Value = 168

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L2").Value = 132
$ws.Range("L4").Value = 28
$ws.Range("L6").Value = 84
$ws.Range("L7").Value = 398

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L2").Value = 115
$ws.Range("L3").Value = 125

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L3").Value = 57
$ws.Range("L7").Value = 193

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("L3").Value = 21
$ws.Range("L7").Value = 79

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L3").Value = 163
$ws.Range("L6").Value = 123
$ws.Range("L7").Value = 497

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("L3").Value = 36

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("L7").Value = 113

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 76

$ws = $wb.Worksheets.Item('East Village')
$ws.Range("L3").Value = 1
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L3").Value = 75
$ws.Range("L7").Value = 241

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("L6").Value = 29

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("L7").Value = 120

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("L6").Value = 65
$ws.Range("L7").Value = 130

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 43

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("L2").Value = 15
$ws.Range("L7").Value = 48

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L6").Value = 162
$ws.Range("L7").Value = 773

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("L2").Value = 9
$ws.Range("L7").Value = 40

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L3").Value = 97
$ws.Range("L4").Value = 19
$ws.Range("L6").Value = 82
$ws.Range("L7").Value = 308

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("L6").Value = 9

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("L7").Value = 35
